# Nalco Aluminium Ingot price sheet — add latest day's quote.
#
# The source feed prepends one new day at the top of the table (row 2) and
# every existing row slides down by one. The row that falls off the bottom
# of the previously-used range (old row 128) survives as the new row 129
# since the table keeps growing.
#
# Strategy: let Excel's own row-insert at row 2 do the heavy lifting (it
# naturally shifts cell values/hyperlinks for rows 3..129 and grows the
# used range/dimension to F129 on its own). Then patch up the two spots
# that a plain insert leaves wrong:
#   - row 2 is blank after the insert: refill B2:F2 from the row that used
#     to be row 2 (now row 3, identical data) and set A2 to the new date.
#   - row 129 has the right values (inherited from old row 128) but no
#     hyperlink yet, since the insert only relocates hyperlinks that were
#     already inside the shifted block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues = -4163

# 1) Push rows 2..128 down to 3..129 (values, styles, and hyperlinks all
#    move together); the blank new row 2 is created with default formatting.
$ws.Range("A2:F2").Insert()

# 2) Give the new row 2 the same look as the data rows (copy from row 3,
#    which now holds what used to be row 2's formatting).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial($xlPasteFormats)

# 3) Row 2's B..F values are identical to the old row 2 (now row 3) — copy
#    them as VALUES (not formats) so the date-looking text in column E
#    stays literal text instead of being reinterpreted as a date.
$ws.Range("B3:F3").Copy()
$ws.Range("B2:F2").PasteSpecial($xlPasteValues)

# 4) Only the date in column A is actually new. Force the cell to Text
#    first so Excel doesn't silently convert "12-12-2025" into a real
#    date serial, then restore the normal data-row formatting (General,
#    centered) on top of the text value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12-12-2025"
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)

# 5) New row 129 (ex-row 128) needs its own hyperlink on F129, matching
#    the link text already sitting in that cell.
$target = $ws.Range("F129").Value()
$ws.Hyperlinks.Add($ws.Range("F129"), $target)

# Re-apply the plain data-row style to F129 — Hyperlinks.Add stamps the
# built-in blue/underline "Hyperlink" style, but every other linked cell
# in this sheet keeps the ordinary centered style.
$ws.Range("F128").Copy()
$ws.Range("F129").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false
